# Applies the cryptos price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "355.80"); force text formatting so Excel
# does not coerce numeric-looking strings into numbers (losing formatting like
# trailing zeros), then clear the format back so no stray cell style is left.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '51.673.50'
$ws.Range("E2").Value = '  -0.35%  '
Set-TextValue $ws.Range("D3") '2.775.61'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  +0.06%  '
Set-TextValue $ws.Range("D5") '355.80'
$ws.Range("E5").Value = '  +0.66%  '
Set-TextValue $ws.Range("D6") '108.74'
$ws.Range("E6").Value = '  -2.14%  '
Set-TextValue $ws.Range("D7") '0.553'
$ws.Range("E7").Value = '  -2.85%  '
$ws.Range("E8").Value = '  +0.09%  '
Set-TextValue $ws.Range("D9") '0.583'
$ws.Range("E9").Value = '  -1.83%  '
Set-TextValue $ws.Range("D10") '39.47'
$ws.Range("E10").Value = '  -2.22%  '
Set-TextValue $ws.Range("D11") '0.136'
$ws.Range("E11").Value = '  +2.99%  '
Set-TextValue $ws.Range("D12") '0.0842'
$ws.Range("E12").Value = '  -1.49%  '
Set-TextValue $ws.Range("D13") '19.36'
$ws.Range("E13").Value = '  -2.50%  '
Set-TextValue $ws.Range("D14") '7.56'
$ws.Range("E14").Value = '  -2.20%  '
Set-TextValue $ws.Range("D15") '3.211.23'
$ws.Range("E15").Value = '  -1.09%  '
Set-TextValue $ws.Range("D16") '2.767.01'
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("E17").Value = '  +0.20%  '
Set-TextValue $ws.Range("D18") '51.649.48'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E19").Value = '  -1.17%  '
Set-TextValue $ws.Range("D20") '3.11'
$ws.Range("E20").Value = '  -0.61%  '
Set-TextValue $ws.Range("D21") '13.05'
$ws.Range("E21").Value = '  -2.29%  '
Set-TextValue $ws.Range("D22") '0.0₃0966'
$ws.Range("E22").Value = '  -2.26%  '
Set-TextValue $ws.Range("D23") '69.99'
$ws.Range("E23").Value = '  -0.61%  '
Set-TextValue $ws.Range("D24") '267.60'
$ws.Range("E24").Value = '  -0.10%  '
Set-TextValue $ws.Range("D25") '2.72'
$ws.Range("E25").Value = '  -2.84%  '
Set-TextValue $ws.Range("D26") '26.31'
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("E27").Value = '  +0.02%  '
Set-TextValue $ws.Range("D28") '0.164'
$ws.Range("E28").Value = '  +16.36%  '
Set-TextValue $ws.Range("D29") '2.30'
$ws.Range("E29").Value = '  +2.01%  '
Set-TextValue $ws.Range("D30") '10.14'
$ws.Range("E30").Value = '  -1.43%  '
Set-TextValue $ws.Range("D31") '6.18'
$ws.Range("E31").Value = '  +4.16%  '
Set-TextValue $ws.Range("D32") '51.50'
$ws.Range("E32").Value = '  -1.94%  '
Set-TextValue $ws.Range("D33") '34.49'
$ws.Range("E33").Value = '  +0.25%  '
Set-TextValue $ws.Range("D34") '0.0449'
$ws.Range("E34").Value = '  -7.34%  '
Set-TextValue $ws.Range("D35") '0.0835'
$ws.Range("E35").Value = '  -1.37%  '
Set-TextValue $ws.Range("D36") '5.13'
$ws.Range("E36").Value = '  -7.27%  '
$ws.Range("E37").Value = '  +0.03%  '
Set-TextValue $ws.Range("D38") '18.68'
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("E40").Value = '  -3.90%  '
$ws.Range("E41").Value = '  +2.46%  '
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("E43").Value = '  -2.52%  '
Set-TextValue $ws.Range("D44") '119.02'
$ws.Range("E44").Value = '  -6.23%  '
Set-TextValue $ws.Range("D45") '21.54'
$ws.Range("E45").Value = '  -6.77%  '
Set-TextValue $ws.Range("D46") '2.081.09'
$ws.Range("E46").Value = '  +0.38%  '
Set-TextValue $ws.Range("D47") '3.24'
$ws.Range("E47").Value = '  -2.20%  '
$ws.Range("E48").Value = '  +1.01%  '
Set-TextValue $ws.Range("D49") '0.930'
$ws.Range("E49").Value = '  -3.36%  '
Set-TextValue $ws.Range("D50") '5.52'
$ws.Range("E50").Value = '  -5.78%  '
Set-TextValue $ws.Range("D51") '8.54'
$ws.Range("E51").Value = '  -6.20%  '
